$wb = $excel.ActiveWorkbook
$cases = $wb.Worksheets.Item("cases")

# Add the brand-new "cases" worksheet right after the existing one so the
# final tab order is generalOptions | cases_doesn't work | cases.
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $cases)
$newSheet.Name = "cases_tmp_new"

# Free up the "cases" name and claim it for the new sheet.
$cases.Name = "cases_doesn't work"
$newSheet.Name = "cases"

# --- re-fetch stable, name-based handles for the rest of the script ---
$oldCases = $wb.Worksheets.Item("cases_doesn't work")
$newCases = $wb.Worksheets.Item("cases")

# note explaining why this sheet doesn't work any more
$oldCases.Range("B7").Value = "this simulation set does not work because there are still bugs with maize appearing on day 256"
$oldCases.Activate()
$oldCases.Range("A1:XFD4").Select()

# populate the new "cases" sheet: wheat + chickpea cases only (maize removed)
$newCases.Range("A1").Value = "name"
$newCases.Range("B1").Value = "climatename"
$newCases.Range("C1").Value = "soilname"
$newCases.Range("D1").Value = "lat"
$newCases.Range("E1").Value = "long"
$newCases.Range("F1").Value = "rotation"
$newCases.Range("G1").Value = "management"

$newCases.Range("A2").Value = "Meknes35degresWheat"
$newCases.Range("B2").Value = "Ain Hamra - Meknes"
$newCases.Range("C2").Value = "325_-35"
$newCases.Range("D2").Value = 35
$newCases.Range("E2").Value = -5
$newCases.Range("F2").Value = '"WHEAT.Ble_Dur_1", "WHEAT.Ble_Tendre_1"'
$newCases.Range("G2").Value = '"ROTATION_BLE", "ROTATION_BLE_IRRIGUE"'

$newCases.Range("A3").Value = "Meknes35degresChickpea"
$newCases.Range("B3").Value = "Ain Hamra - Meknes"
$newCases.Range("C3").Value = "325_-35"
$newCases.Range("D3").Value = 45
$newCases.Range("E3").Value = -5
$newCases.Range("F3").Value = '"Chickpea.Ghab2", "WHEAT.Cocorit", "WHEAT.Avoine_Romani"'
$newCases.Range("G3").Value = '"ROTATION_POISCHICHE", "ROTATION_BLE", "ROTATION_BLE_IRRIGUE"'

$newCases.Activate()
$newCases.Range("F9").Select()

$newCases.PageSetup.LeftMargin = 54
$newCases.PageSetup.RightMargin = 54
$newCases.PageSetup.TopMargin = 72
$newCases.PageSetup.BottomMargin = 72
$newCases.PageSetup.HeaderMargin = 36
$newCases.PageSetup.FooterMargin = 36

Write-Output "All sheets in order:"
foreach ($s in $wb.Worksheets) { Write-Output $s.Name }
